$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Agrícola del Norte S.A. de Arica - Tuna
$ws.Range("D2").Value = 44679
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 29000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 29500
$ws.Range("Q2").Value = "$/caja 20 kilos"
$ws.Range("S2").Value = 1475
$ws.Range("T2").Value = 20

# Row 3
$ws.Range("D3").Value = 44679
$ws.Range("L3").Value = "Tercera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("Q3").Value = "$/caja 20 kilos"
$ws.Range("S3").Value = 1225
$ws.Range("T3").Value = 20

# Row 4
$ws.Range("D4").Value = 44643
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 28000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 29000
$ws.Range("Q4").Value = "$/caja 20 kilos"
$ws.Range("S4").Value = 1450
$ws.Range("T4").Value = 20

# Row 5
$ws.Range("D5").Value = 44664
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 29000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 29500
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("S5").Value = 1639
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44650
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 31000
$ws.Range("O6").Value = 32000
$ws.Range("P6").Value = 31500
$ws.Range("Q6").Value = "$/caja 20 kilos"
$ws.Range("S6").Value = 1575
$ws.Range("T6").Value = 20

# Row 7
$ws.Range("D7").Value = 44650
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("Q7").Value = "$/caja 20 kilos"
$ws.Range("S7").Value = 1475
$ws.Range("T7").Value = 20

# Row 8
$ws.Range("D8").Value = 44671
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 29000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 29500
$ws.Range("Q8").Value = "$/caja 20 kilos"
$ws.Range("S8").Value = 1475
$ws.Range("T8").Value = 20

# Row 9
$ws.Range("D9").Value = 44636
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 29000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 29500
$ws.Range("Q9").Value = "$/caja 20 kilos"
$ws.Range("S9").Value = 1475
$ws.Range("T9").Value = 20
